$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated values for column F (dSF) on specific rows (repull data / recalculation)
$updates = @{
    11 = 1
    15 = 3
    16 = 3
    17 = 1
    19 = 0
    34 = -1
    35 = 1
    43 = -2
    48 = 4
    51 = 5
    53 = -2
    61 = 0
    62 = -3
    69 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
